# Applies the crypto-price/volume refresh described by the commit:
# "Updated cryptos list on Sat May  6 09:11:32 UTC 2023 with GitHub Actions"
#
# Column D ("Price") cells are stored as literal text in the workbook (e.g. "46.17",
# "29.418.94") rather than numbers, so that values such as trailing zeros and
# thousand-grouped figures like "29.418.94" render exactly as scraped. Any new price
# whose text happens to look like a plain number (e.g. "46.10") would otherwise be
# auto-converted to a Number by Excel on assignment (losing the trailing zero), so
# those cells are pre-formatted as Text ("@") before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "Price" (column D) values that look like plain numbers: force Text format
# first so Excel keeps the exact scraped string instead of coercing to a Number. ---
$priceTextUpdates = [ordered]@{
    'D5' = '325.81'
    'D6' = '1.001'
    'D8' = '0.3871'
    'D9' = '46.10'
    'D10' = '0.07833'
    'D11' = '0.9768'
    'D12' = '22.66'
    'D14' = '7.085'
    'D15' = '5.757'
    'D16' = '0.07049'
    'D17' = '86.78'
    'D19' = '0.000009820'
    'D20' = '17.08'
    'D21' = '1.001'
    'D23' = '5.472'
    'D24' = '11.09'
    'D26' = '2.097'
    'D27' = '157.28'
    'D28' = '19.39'
    'D29' = '5.769'
    'D30' = '118.62'
    'D31' = '1.868'
    'D32' = '0.09373'
    'D33' = '0.8641'
    'D34' = '5.195'
    'D35' = '1.308'
    'D36' = '3.124'
    'D37' = '0.05773'
    'D38' = '1.154'
    'D39' = '0.02086'
    'D40' = '7.721'
    'D41' = '0.5664'
    'D42' = '0.1783'
    'D43' = '9.442'
    'D44' = '0.000002903'
    'D45' = '2.724'
    'D46' = '11.63'
    'D47' = '0.5296'
    'D48' = '2.092'
    'D49' = '0.06876'
    'D50' = '1.819'
    'D51' = '111.46'
}
foreach ($cell in $priceTextUpdates.Keys) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $priceTextUpdates[$cell]
}

# --- All other updated cells (Coin / Link / non-numeric Price / Volume%) ---
$otherUpdates = [ordered]@{
    'D2' = '29.416.09'
    'E2' = '  +1.10%  '
    'D3' = '1.947.06'
    'E3' = '  +2.53%  '
    'E4' = '  +0.15%  '
    'E5' = '  +0.23%  '
    'E6' = '  +0.10%  '
    'E7' = '  +0.77%  '
    'E8' = '  -0.30%  '
    'E9' = '  +0.50%  '
    'E10' = '  -0.37%  '
    'E11' = '  -1.22%  '
    'E12' = '  +3.65%  '
    'D13' = '1.935.97'
    'E13' = '  +2.20%  '
    'E14' = '  +0.52%  '
    'E15' = '  -0.28%  '
    'E16' = '  +0.59%  '
    'E17' = '  -1.22%  '
    'E19' = '  -1.17%  '
    'E20' = '  +0.31%  '
    'E21' = '  +0.11%  '
    'D22' = '29.414.66'
    'E22' = '  +1.08%  '
    'E23' = '  +2.85%  '
    'D25' = '2.169.82'
    'E25' = '  +2.03%  '
    'E26' = '  -0.11%  '
    'E27' = '  +0.71%  '
    'E28' = '  +0.09%  '
    'E29' = '  -2.29%  '
    'E30' = '  +0.22%  '
    'E31' = '  +0.04%  '
    'E32' = '  +0.65%  '
    'E33' = '  -3.52%  '
    'E34' = '  -0.72%  '
    'E35' = '  -0.72%  '
    'E36' = '  -0.51%  '
    'E37' = '  -0.11%  '
    'E38' = '  -1.17%  '
    'E39' = '  +0.13%  '
    'E40' = '  +0.67%  '
    'E41' = '  -0.15%  '
    'E42' = '  -0.55%  '
    'E43' = '  -2.52%  '
    'E44' = '  +37.74%  '
    'E45' = '  +6.77%  '
    'B46' = 'EnergySwap'
    'C46' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E46' = '  -1.96%  '
    'B47' = 'Decentraland'
    'C47' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'E47' = '  -0.98%  '
    'E48' = '  -5.55%  '
    'E49' = '  -1.82%  '
    'E50' = '  -1.45%  '
    'E51' = '  -0.99%  '
}
foreach ($cell in $otherUpdates.Keys) {
    $ws.Range($cell).Value = $otherUpdates[$cell]
}
